# Auto-generated script applying value updates to Coeurl_Profits workbook
# Data refresh from scheduled runner: updates computed market-price columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all class sheets.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3175
$ws.Range("I51").Value = 2750.0625
$ws.Range("K51").Value = 2750.0625
$ws.Range("M51").Value = -2266.0625
$ws.Range("H113").Value = 103480
$ws.Range("I113").Value = 179000
$ws.Range("J113").Value = 12856
$ws.Range("K113").Value = 179000
$ws.Range("L113").Value = 12856
$ws.Range("M113").Value = -175746
$ws.Range("N113").Value = -19364
$ws.Range("H132").Value = 42236.285
$ws.Range("I132").Value = 24534.488
$ws.Range("J132").Value = 169099.17
$ws.Range("K132").Value = 73603.46400000001
$ws.Range("L132").Value = 507297.51
$ws.Range("M132").Value = -71073.46400000001
$ws.Range("N132").Value = -512357.51
$ws.Range("H136").Value = 49874.875
$ws.Range("J136").Value = 49874.875
$ws.Range("L136").Value = 49874.875
$ws.Range("N136").Value = -60074.875
$ws.Range("H137").Value = 1526787.4
$ws.Range("I137").Value = 43583.46
$ws.Range("K137").Value = 130750.38
$ws.Range("M137").Value = -128200.38
$ws.Range("H138").Value = 4052.642
$ws.Range("J138").Value = 4523.179
$ws.Range("L138").Value = 13569.537
$ws.Range("N138").Value = -23849.537

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1528.6154
$ws.Range("I2").Value = 1447.6666
$ws.Range("K2").Value = 1447.6666
$ws.Range("M2").Value = -1334.6666
$ws.Range("H63").Value = 3426.6667
$ws.Range("I63").Value = 3426.6667
$ws.Range("K63").Value = 3426.6667
$ws.Range("M63").Value = -2740.6667
$ws.Range("H66").Value = 3426.6667
$ws.Range("I66").Value = 3426.6667
$ws.Range("K66").Value = 17133.3335
$ws.Range("M66").Value = -13701.3335
$ws.Range("H80").Value = 42055
$ws.Range("I80").Value = 10000
$ws.Range("K80").Value = 10000
$ws.Range("M80").Value = -9002
$ws.Range("H83").Value = 42055
$ws.Range("I83").Value = 10000
$ws.Range("K83").Value = 30000
$ws.Range("M83").Value = -25008
$ws.Range("H116").Value = 1528.6154
$ws.Range("I116").Value = 1447.6666
$ws.Range("K116").Value = 1447.6666
$ws.Range("M116").Value = 846.3334
$ws.Range("H132").Value = 3459.5518
$ws.Range("I132").Value = 3124.7896
$ws.Range("J132").Value = 4095.6
$ws.Range("K132").Value = 9374.3688
$ws.Range("L132").Value = 12286.8
$ws.Range("M132").Value = -6844.3688
$ws.Range("N132").Value = -17346.8

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1528.6154
$ws.Range("I3").Value = 1447.6666
$ws.Range("K3").Value = 1447.6666
$ws.Range("M3").Value = -1333.6666
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H86").Value = 4031.2258
$ws.Range("I86").Value = 3096.4285
$ws.Range("J86").Value = 5994.3
$ws.Range("K86").Value = 3096.4285
$ws.Range("L86").Value = 5994.3
$ws.Range("M86").Value = -1973.4285
$ws.Range("N86").Value = -8240.299999999999
$ws.Range("H89").Value = 4031.2258
$ws.Range("I89").Value = 3096.4285
$ws.Range("J89").Value = 5994.3
$ws.Range("K89").Value = 15482.1425
$ws.Range("L89").Value = 29971.5
$ws.Range("M89").Value = -9866.1425
$ws.Range("N89").Value = -41203.5
$ws.Range("H94").Value = 1526.7778
$ws.Range("I94").Value = 1578.1428
$ws.Range("J94").Value = 1347
$ws.Range("K94").Value = 1578.1428
$ws.Range("L94").Value = 1347
$ws.Range("M94").Value = -1127.1428
$ws.Range("N94").Value = -2249
$ws.Range("H105").Value = 8666.296
$ws.Range("I105").Value = 7847.391
$ws.Range("K105").Value = 7847.391
$ws.Range("M105").Value = -6100.391
$ws.Range("H140").Value = 71993.75
$ws.Range("J140").Value = 71993.75
$ws.Range("L140").Value = 71993.75
$ws.Range("N140").Value = -82353.75

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 205.16667
$ws.Range("I7").Value = 109
$ws.Range("J7").Value = 356.2857
$ws.Range("K7").Value = 109
$ws.Range("L7").Value = 356.2857
$ws.Range("M7").Value = 4
$ws.Range("N7").Value = -582.2857
$ws.Range("H62").Value = 5177.2383
$ws.Range("I62").Value = 4734.8887
$ws.Range("J62").Value = 7831.3335
$ws.Range("K62").Value = 4734.8887
$ws.Range("L62").Value = 7831.3335
$ws.Range("M62").Value = -4110.8887
$ws.Range("N62").Value = -9079.333500000001
$ws.Range("H65").Value = 5177.2383
$ws.Range("I65").Value = 4734.8887
$ws.Range("J65").Value = 7831.3335
$ws.Range("K65").Value = 23674.4435
$ws.Range("L65").Value = 39156.6675
$ws.Range("M65").Value = -20554.4435
$ws.Range("N65").Value = -45396.6675
$ws.Range("H105").Value = 3425
$ws.Range("I105").Value = 3425
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3425
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1678
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 2155.6304
$ws.Range("I132").Value = 1690.1892
$ws.Range("K132").Value = 5070.5676
$ws.Range("M132").Value = -2540.5676

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 22222570
$ws.Range("I60").Value = 33333704
$ws.Range("J60").Value = 299.2
$ws.Range("K60").Value = 100001112
$ws.Range("L60").Value = 897.5999999999999
$ws.Range("M60").Value = -100000861
$ws.Range("N60").Value = -1399.6
$ws.Range("H68").Value = 1924651.4
$ws.Range("J68").Value = 2274419.8
$ws.Range("L68").Value = 6823259.399999999
$ws.Range("N68").Value = -6824881.399999999
$ws.Range("H71").Value = 1924651.4
$ws.Range("J71").Value = 2274419.8
$ws.Range("L71").Value = 20469778.2
$ws.Range("N71").Value = -20477890.2
$ws.Range("H124").Value = 9546.718999999999
$ws.Range("J124").Value = 9796.5
$ws.Range("L124").Value = 29389.5
$ws.Range("N124").Value = -39209.5
$ws.Range("H129").Value = 2189.889
$ws.Range("I129").Value = 762.75
$ws.Range("K129").Value = 2288.25
$ws.Range("M129").Value = 2711.75
$ws.Range("H137").Value = 4852.615
$ws.Range("J137").Value = 2485
$ws.Range("L137").Value = 7455
$ws.Range("N137").Value = -17655

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2452.647
$ws.Range("I80").Value = 2478.2144
$ws.Range("J80").Value = 2333.3333
$ws.Range("K80").Value = 2478.2144
$ws.Range("L80").Value = 2333.3333
$ws.Range("M80").Value = -1480.2144
$ws.Range("N80").Value = -4329.3333
$ws.Range("H83").Value = 2452.647
$ws.Range("I83").Value = 2478.2144
$ws.Range("J83").Value = 2333.3333
$ws.Range("K83").Value = 12391.072
$ws.Range("L83").Value = 11666.6665
$ws.Range("M83").Value = -7399.072
$ws.Range("N83").Value = -21650.6665
$ws.Range("H102").Value = 1541.4
$ws.Range("I102").Value = 1476.75
$ws.Range("J102").Value = 1800
$ws.Range("K102").Value = 1476.75
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = 145.25
$ws.Range("N102").Value = -5044
$ws.Range("H113").Value = 4799.478
$ws.Range("I113").Value = 4262.8184
$ws.Range("K113").Value = 4262.8184
$ws.Range("M113").Value = -2092.8184
$ws.Range("H122").Value = 5707.6
$ws.Range("I122").Value = 5803.1333
$ws.Range("J122").Value = 5421
$ws.Range("K122").Value = 17409.3999
$ws.Range("L122").Value = 16263
$ws.Range("M122").Value = -14959.3999
$ws.Range("N122").Value = -21163
$ws.Range("H132").Value = 14513.525
$ws.Range("J132").Value = 3057
$ws.Range("L132").Value = 9171
$ws.Range("N132").Value = -14231
$ws.Range("H135").Value = 49500
$ws.Range("J135").Value = 49500
$ws.Range("L135").Value = 49500
$ws.Range("N135").Value = -59640

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5288.2
$ws.Range("I16").Value = 5556
$ws.Range("J16").Value = 4663.3335
$ws.Range("K16").Value = 5556
$ws.Range("L16").Value = 4663.3335
$ws.Range("M16").Value = -5386
$ws.Range("N16").Value = -5003.3335
$ws.Range("H40").Value = 7069.7646
$ws.Range("I40").Value = 5854.5557
$ws.Range("K40").Value = 5854.5557
$ws.Range("M40").Value = -5718.5557
$ws.Range("H46").Value = 2138.7646
$ws.Range("J46").Value = 3147.125
$ws.Range("L46").Value = 3147.125
$ws.Range("N46").Value = -3523.125
$ws.Range("H122").Value = 6174
$ws.Range("I122").Value = 5703
$ws.Range("K122").Value = 17109
$ws.Range("M122").Value = -14659

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19227.912
$ws.Range("I62").Value = 18630.691
$ws.Range("K62").Value = 18630.691
$ws.Range("M62").Value = -18006.691
$ws.Range("H65").Value = 19227.912
$ws.Range("I65").Value = 18630.691
$ws.Range("K65").Value = 93153.45499999999
$ws.Range("M65").Value = -90033.45499999999
$ws.Range("H81").Value = 7598.619
$ws.Range("J81").Value = 3815.7144
$ws.Range("L81").Value = 7631.4288
$ws.Range("N81").Value = -9753.4288
$ws.Range("H84").Value = 7598.619
$ws.Range("J84").Value = 3815.7144
$ws.Range("L84").Value = 38157.144
$ws.Range("N84").Value = -48765.144
$ws.Range("H122").Value = 4951.4585
$ws.Range("I122").Value = 3055.0625
$ws.Range("J122").Value = 8744.25
$ws.Range("K122").Value = 9165.1875
$ws.Range("L122").Value = 26232.75
$ws.Range("M122").Value = -6715.1875
$ws.Range("N122").Value = -31132.75
$ws.Range("H136").Value = 1587
$ws.Range("I136").Value = 1202.8636
$ws.Range("J136").Value = 2995.5
$ws.Range("K136").Value = 3608.5908
$ws.Range("L136").Value = 8986.5
$ws.Range("M136").Value = -1058.5908
$ws.Range("N136").Value = -14086.5
